$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 5.392984666666666
$ws.Range("H2").Value = 16.178954
$ws.Range("I2").Value = 0.3709566379599203
$ws.Range("J2").Value = 0.3709566379599202
$ws.Range("K2").Value = 2
$ws.Range("M2").Value = 86.777428
$ws.Range("N2").Value = 173.554856
$ws.Range("O2").Value = 0.2676090626666408
$ws.Range("P2").Value = 0.2041884050300022
$ws.Range("Q2").Value = 467.9893386167706
$ws.Range("R2").Value = 2807.936031700624
$ws.Range("S2").Value = 0.09927135817442269
$ws.Range("T2").Value = 0.07574504424032809

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 5.392984666666666
$ws.Range("H3").Value = 16.178954
$ws.Range("I3").Value = 0.3709566379599203
$ws.Range("J3").Value = 0.3709566379599202
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 128.9086913333333
$ws.Range("N3").Value = 386.726074
$ws.Range("O3").Value = 0.3975357976419474
$ws.Range("P3").Value = 0.4549857149118007
$ws.Range("Q3").Value = 695.2025957607328
$ws.Range("R3").Value = 6256.823361846596
$ws.Range("S3").Value = 0.147468542961972
$ws.Range("T3").Value = 0.1687799711234723

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 5.392984666666666
$ws.Range("H4").Value = 16.178954
$ws.Range("I4").Value = 0.3709566379599203
$ws.Range("J4").Value = 0.3709566379599202
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 12.52958833333333
$ws.Range("N4").Value = 37.588765
$ws.Range("O4").Value = 0.03863944192356349
$ws.Range("P4").Value = 0.04422342393230168
$ws.Range("Q4").Value = 67.57187776131222
$ws.Range("R4").Value = 608.14689985181
$ws.Range("S4").Value = 0.01433355746861271
$ws.Range("T4").Value = 0.01640497266100291

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 5.392984666666666
$ws.Range("H5").Value = 16.178954
$ws.Range("I5").Value = 0.3709566379599203
$ws.Range("J5").Value = 0.3709566379599202
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 17.279662
$ws.Range("N5").Value = 51.83898599999999
$ws.Range("O5").Value = 0.05328798349515926
$ws.Range("P5").Value = 0.06098890064886812
$ws.Range("Q5").Value = 93.18895221118265
$ws.Range("R5").Value = 838.7005699006437
$ws.Range("S5").Value = 0.019767531201028
$ws.Range("T5").Value = 0.02262423753757572

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 5.392984666666666
$ws.Range("H6").Value = 16.178954
$ws.Range("I6").Value = 0.3709566379599203
$ws.Range("J6").Value = 0.3709566379599202
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 42.71737566666666
$ws.Range("N6").Value = 128.152127
$ws.Range("O6").Value = 0.1317342208129911
$ws.Range("P6").Value = 0.1507718021634167
$ws.Range("Q6").Value = 230.3741519705731
$ws.Range("R6").Value = 2073.367367735157
$ws.Range("S6").Value = 0.04886768365705694
$ws.Range("T6").Value = 0.0559298008296993

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 5.392984666666666
$ws.Range("H7").Value = 16.178954
$ws.Range("I7").Value = 0.3709566379599203
$ws.Range("J7").Value = 0.3709566379599202
$ws.Range("K7").Value = 2
$ws.Range("M7").Value = 36.0566465
$ws.Range("N7").Value = 72.113293
$ws.Range("O7").Value = 0.111193493459698
$ws.Range("P7").Value = 0.08484175331361067
$ws.Range("Q7").Value = 194.4529417059203
$ws.Range("R7").Value = 1166.717650235522
$ws.Range("S7").Value = 0.04124796449682796
$ws.Range("T7").Value = 0.03147261156784194

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 9.145061333333333
$ws.Range("H8").Value = 27.435184
$ws.Range("I8").Value = 0.6290433620400798
$ws.Range("J8").Value = 0.6290433620400798
$ws.Range("K8").Value = 2
$ws.Range("M8").Value = 86.777428
$ws.Range("N8").Value = 173.554856
$ws.Range("O8").Value = 0.2676090626666408
$ws.Range("P8").Value = 0.2041884050300022
$ws.Range("Q8").Value = 793.5849014089173
$ws.Range("R8").Value = 4761.509408453504
$ws.Range("S8").Value = 0.1683377044922181
$ws.Range("T8").Value = 0.1284433607896741

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 9.145061333333333
$ws.Range("H9").Value = 27.435184
$ws.Range("I9").Value = 0.6290433620400798
$ws.Range("J9").Value = 0.6290433620400798
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 128.9086913333333
$ws.Range("N9").Value = 386.726074
$ws.Range("O9").Value = 0.3975357976419474
$ws.Range("P9").Value = 0.4549857149118007
$ws.Range("Q9").Value = 1178.877888643068
$ws.Range("R9").Value = 10609.90099778762
$ws.Range("S9").Value = 0.2500672546799754
$ws.Range("T9").Value = 0.2862057437883284

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 9.145061333333333
$ws.Range("H10").Value = 27.435184
$ws.Range("I10").Value = 0.6290433620400798
$ws.Range("J10").Value = 0.6290433620400798
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 12.52958833333333
$ws.Range("N10").Value = 37.588765
$ws.Range("O10").Value = 0.03863944192356349
$ws.Range("P10").Value = 0.04422342393230168
$ws.Range("Q10").Value = 114.5838537897511
$ws.Range("R10").Value = 1031.25468410776
$ws.Range("S10").Value = 0.02430588445495079
$ws.Range("T10").Value = 0.02781845127129878

$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 9.145061333333333
$ws.Range("H11").Value = 27.435184
$ws.Range("I11").Value = 0.6290433620400798
$ws.Range("J11").Value = 0.6290433620400798
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 17.279662
$ws.Range("N11").Value = 51.83898599999999
$ws.Range("O11").Value = 0.05328798349515926
$ws.Range("P11").Value = 0.06098890064886812
$ws.Range("Q11").Value = 158.0235688092693
$ws.Range("R11").Value = 1422.212119283424
$ws.Range("S11").Value = 0.03352045229413127
$ws.Range("T11").Value = 0.03836466311129241

$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 9.145061333333333
$ws.Range("H12").Value = 27.435184
$ws.Range("I12").Value = 0.6290433620400798
$ws.Range("J12").Value = 0.6290433620400798
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 42.71737566666666
$ws.Range("N12").Value = 128.152127
$ws.Range("O12").Value = 0.1317342208129911
$ws.Range("P12").Value = 0.1507718021634167
$ws.Range("Q12").Value = 390.6530204707075
$ws.Range("R12").Value = 3515.877184236367
$ws.Range("S12").Value = 0.08286653715593419
$ws.Range("T12").Value = 0.09484200133371744

$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 9.145061333333333
$ws.Range("H13").Value = 27.435184
$ws.Range("I13").Value = 0.6290433620400798
$ws.Range("J13").Value = 0.6290433620400798
$ws.Range("K13").Value = 2
$ws.Range("M13").Value = 36.0566465
$ws.Range("N13").Value = 72.113293
$ws.Range("O13").Value = 0.111193493459698
$ws.Range("P13").Value = 0.08484175331361067
$ws.Range("Q13").Value = 329.7402437168186
$ws.Range("R13").Value = 1978.441462300912
$ws.Range("S13").Value = 0.06994552896287005
$ws.Range("T13").Value = 0.05336914174576874

